# Update the LR-pairs data for Fgf8-Fgfr4 with new TPM values, and drop
# the rows for the cell-type pairs that are no longer included.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 7-9 (the sheet used to run through row 9; now it only
# needs rows through row 6).
$ws.Range("A7:T9").EntireRow.Delete() | Out-Null

# Target (cluster) names for the remaining 5 data rows (2..6).
$targets = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")

# Numeric data (columns E..T) for each of the 5 data rows.
$data = @(
    @(1, 0.3333333333333333, 0.1008186666666667, 0.302456, 1, 1, 3, 1, 0.363908, 1.091724, 0.01118972054940699, 0.01663265480083782, 0.03668871934933333, 0.330198474144, 0.01118972054940699, 0.01663265480083782),
    @(1, 0.3333333333333333, 0.1008186666666667, 0.302456, 1, 1, 1, 0.3333333333333333, 0.153566, 0.460698, 0.004721964413781051, 0.007018835164781924, 0.01548231936533333, 0.139340874288, 0.004721964413781051, 0.007018835164781924),
    @(1, 0.3333333333333333, 0.1008186666666667, 0.302456, 1, 1, 1, 0.3333333333333333, 0.06252866666666666, 0.187586, 0.001922679101110775, 0.002857913889838424, 0.006304056801777778, 0.056736511216, 0.001922679101110775, 0.002857913889838424),
    @(1, 0.3333333333333333, 0.1008186666666667, 0.302456, 1, 1, 2, 1, 31.927516, 63.85503199999999, 0.9817315966582778, 0.9728454303033116, 3.218889593098667, 19.313337558592, 0.9817315966582778, 0.9728454303033116),
    @(1, 0.3333333333333333, 0.1008186666666667, 0.302456, 1, 1, 1, 0.3333333333333333, 0.01411566666666667, 0.042347, 0.0004340392774233579, 0.0006451658412300904, 0.001423122692444444, 0.012808104232, 0.0004340392774233579, 0.0006451658412300904)
)

for ($i = 0; $i -lt $targets.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "ECs"
    $ws.Cells.Item($row, 2).Value = "Fgf8"
    $ws.Cells.Item($row, 3).Value = "Fgfr4"
    $ws.Cells.Item($row, 4).Value = $targets[$i]

    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 5 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
